# regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals — the recomputed "K" (column G) values per game row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1;  3  = 3;  4  = 3;  5  = 3;  6  = 2;  7  = 2;  8  = 2;  9  = 0
    10 = 1;  11 = 1;  12 = 2;  13 = 1;  14 = 0;  15 = 0;  16 = 3;  17 = 1
    18 = 1;  19 = 1;  20 = 1;  21 = 2;  22 = 0;  24 = 2;  25 = 1;  26 = 1
    27 = 3;  28 = 1;  29 = 2;  30 = 1;  31 = 1;  32 = 2;  33 = 3;  34 = 3
    35 = 0;  36 = 1;  37 = 1;  38 = 1;  39 = 2;  40 = 2;  41 = 0;  42 = 1
    43 = 0;  44 = 3;  45 = 1;  46 = 1;  47 = 0;  48 = 2;  49 = 2;  50 = 2
    51 = 3;  52 = 3;  53 = 2;  54 = 1;  55 = 1;  56 = 1;  57 = 3;  58 = 1
    59 = 1;  60 = 1;  61 = 1;  62 = 1;  63 = 2;  64 = 1;  65 = 3;  66 = 2
    67 = 1;  68 = 2;  69 = 1;  70 = 3;  71 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
